$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "stretchy taxi" / "com.singleton.strechy" row (row 2) was removed,
# shifting every subsequent row up by one.
$ws.Rows(2).Delete()

# Two rows were appended at the bottom: "taxi game" / "com.singleton.strechy"
# and "taxi" / "com.singleton.strechy". Seed them from the last existing
# row's formatting (copy/paste-special) so they pick up the same cell style
# instead of Excel's default style.
$ws.Range("A18:B18").Copy()
$ws.Range("A19:B20").PasteSpecial(-4122)

$ws.Range("A19").Value = "taxi game"
$ws.Range("B19").Value = "com.singleton.strechy"
$ws.Range("A20").Value = "taxi"
$ws.Range("B20").Value = "com.singleton.strechy"
